$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 ----
# Preserve the existing D3/E3/F3 values by copying them over to I3/J3/K3
# (Copy/paste keeps the original cell data type - number vs. shared string -
# intact, unlike re-assigning .Value2 which would coerce numeric-looking
# text back into a real number).
$ws.Range("D3").Copy($ws.Range("I3"))
$ws.Range("E3").Copy($ws.Range("J3"))
$ws.Range("F3").Copy($ws.Range("K3"))
# Write the new, updated D3/E3/F3 values
$ws.Range("D3").Value2 = 1.8
$ws.Range("E3").Value2 = "6.41, 9.46"
$ws.Range("F3").Value2 = 0.2

# ---- Row 4 ----
$ws.Range("D4").Copy($ws.Range("I4"))
$ws.Range("E4").Copy($ws.Range("J4"))
$ws.Range("F4").Copy($ws.Range("K4"))
$ws.Range("D4").Value2 = 0.87
$ws.Range("E4").Value2 = "3.23, 4.72"
$ws.Range("F4").Value2 = 0.08

# ---- Row 5 ----
$ws.Range("D5").Copy($ws.Range("I5"))
$ws.Range("E5").Copy($ws.Range("J5"))
$ws.Range("F5").Copy($ws.Range("K5"))
$ws.Range("D5").Value2 = 0.49
$ws.Range("E5").Value2 = "2.6, 3.37"
$ws.Range("F5").Value2 = 0.04

# Update the active selection to match the saved workbook state
$ws.Range("F5").Select()
